$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.140.02"
$ws.Range("E2").Value = "  -0.05%  "
$ws.Range("D3").Value = "1.832.96"
$ws.Range("E3").Value = "  -0.32%  "
$ws.Range("D4").Value = "'0.9993"
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "'241.32"
$ws.Range("E5").Value = "  +0.55%  "
$ws.Range("D6").Value = "'0.6642"
$ws.Range("E6").Value = "  -2.61%  "
$ws.Range("D7").Value = "'0.9999"
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("D8").Value = "'0.07427"
$ws.Range("E8").Value = "  -0.34%  "
$ws.Range("E9").Value = "  -1.87%  "
$ws.Range("D10").Value = "'22.72"
$ws.Range("E10").Value = "  -2.18%  "
$ws.Range("D11").Value = "'0.07736"
$ws.Range("E11").Value = "  +1.20%  "
$ws.Range("D12").Value = "1.840.48"
$ws.Range("E12").Value = "  +0.43%  "
$ws.Range("E13").Value = "  -0.76%  "
$ws.Range("D14").Value = "'0.6681"
$ws.Range("E14").Value = "  -1.89%  "
$ws.Range("D15").Value = "'82.75"
$ws.Range("E15").Value = "  -5.22%  "
$ws.Range("D16").Value = "'6.097"
$ws.Range("E16").Value = "  -0.96%  "
$ws.Range("D17").Value = "'0.000008362"
$ws.Range("E17").Value = "  +1.88%  "
$ws.Range("D18").Value = "29.136.59"
$ws.Range("E18").Value = "  +0.02%  "
$ws.Range("D19").Value = "'227.48"
$ws.Range("E19").Value = "  -1.09%  "
$ws.Range("D20").Value = "'12.48"
$ws.Range("E20").Value = "  -0.24%  "
$ws.Range("D22").Value = "'7.158"
$ws.Range("E22").Value = "  -2.60%  "
$ws.Range("E23").Value = "  +0.04%  "
$ws.Range("D24").Value = "'159.80"
$ws.Range("E24").Value = "  -0.72%  "
$ws.Range("E25").Value = "  -1.09%  "
$ws.Range("D26").Value = "'0.1402"
$ws.Range("E26").Value = "  -1.90%  "
$ws.Range("D27").Value = "'17.96"
$ws.Range("E27").Value = "  -0.53%  "
$ws.Range("D28").Value = "'1.508"
$ws.Range("E28").Value = "  +0.24%  "
$ws.Range("D29").Value = "'4.114"
$ws.Range("E29").Value = "  -3.26%  "
$ws.Range("E30").Value = "  -2.37%  "
$ws.Range("E31").Value = "  -0.10%  "
$ws.Range("D32").Value = "'0.05311"
$ws.Range("E32").Value = "  -0.53%  "
$ws.Range("E33").Value = "  +1.06%  "
$ws.Range("D34").Value = "'0.7518"
$ws.Range("E34").Value = "  -0.35%  "
$ws.Range("E35").Value = "  +0.34%  "
$ws.Range("D36").Value = "'2.652"
$ws.Range("E36").Value = "  -1.10%  "
$ws.Range("E37").Value = "  -2.41%  "
$ws.Range("D38").Value = "'0.01795"
$ws.Range("D39").Value = "'2.732"
$ws.Range("E39").Value = "  +0.53%  "
$ws.Range("D40").Value = "'0.9286"
$ws.Range("E40").Value = "  -2.07%  "
$ws.Range("D41").Value = "'0.08715"
$ws.Range("E41").Value = "  +6.97%  "
$ws.Range("D42").Value = "'5.931"
$ws.Range("E42").Value = "  -2.19%  "
$ws.Range("D43").Value = "'0.9998"
$ws.Range("E43").Value = "  +0.09%  "
$ws.Range("D44").Value = "'101.81"
$ws.Range("E44").Value = "  -3.21%  "
$ws.Range("D45").Value = "1.976.85"
$ws.Range("E45").Value = "  +0.12%  "
$ws.Range("D46").Value = "'0.5144"
$ws.Range("E46").Value = "  -0.60%  "
$ws.Range("E47").Value = "  -0.53%  "
$ws.Range("E48").Value = "  -1.38%  "
$ws.Range("D49").Value = "'63.35"
$ws.Range("E49").Value = "  -1.23%  "
$ws.Range("D50").Value = "'0.05887"
$ws.Range("E50").Value = "  -0.94%  "
$ws.Range("D51").Value = "'6.789"
$ws.Range("E51").Value = "  -1.51%  "
